# Trade #52 closed at 2026-02-17 08:42:12 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and "MarketMaking"
# sheets of the live trading results workbook to record the newly closed
# trade (#52, index 52 in the log) on the MarketMaking strategy.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Summary" - aggregate statistics across all strategies
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.44            # Current Capital
$summary.Range("B4").Value = -0.5600000000000001 # Total P&L $
$summary.Range("B5").Value = -0.22              # Total P&L %
$summary.Range("B6").Value = 52                 # Total Trades
$summary.Range("B8").Value = 24                 # Losing Trades
$summary.Range("B9").Value = 34.62              # Win Rate %

# ---------------------------------------------------------------------
# Sheet 2: "Strategy Status" - per-strategy snapshot (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.44                # Capital
$status.Range("D4").Value = 52                   # Trades
$status.Range("E4").Value = -0.5600000000000001  # P&L $
$status.Range("F4").Value = -0.5600000000000001  # P&L %
$status.Range("G4").Value = 34.62                # Win Rate %

# ---------------------------------------------------------------------
# Helper: append the new trade row (#52) to a trade log worksheet
# ---------------------------------------------------------------------
function Add-TradeRow52($ws) {
    $row = 53

    $ws.Cells.Item($row, 1).Value = 52                # Trade #

    # Use a leading apostrophe so Excel stores the date/time as literal
    # text instead of auto-converting it to a date serial number, which
    # matches how the rest of the sheet stores these values.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"     # Date
    $ws.Cells.Item($row, 3).Value = "08:42:06"        # Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"    # Strategy
    $ws.Cells.Item($row, 5).Value = "UP"              # Side
    $ws.Cells.Item($row, 6).Value = 0.61              # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.59              # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"          # Status
    $ws.Cells.Item($row, 9).Value = -3.2787           # P&L %
    $ws.Cells.Item($row, 10).Value = -0.02            # P&L $
    $ws.Cells.Item($row, 11).Value = 99.44            # Capital After
    $ws.Cells.Item($row, 12).Value = 0                # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"     # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.12             # Duration (min)
}

# ---------------------------------------------------------------------
# Sheet 3: "All Trades" - full trade log
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow52 $allTrades

# ---------------------------------------------------------------------
# Sheet 4: "MarketMaking" - strategy-specific trade log (same content)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow52 $marketMaking
